$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores pre-formatted text such as "166.36" or
# "68.324.09" rather than real numbers. Force the cell to text format
# before writing so Excel does not reinterpret it as a number (which would
# drop trailing zeros / merge the thousands separators), then restore the
# default "Normal" style so no new formatting is left behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.388.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.746.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.746.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.373.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.753.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.404.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.701.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "391.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.748.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.78%  "
